# Updates cryptos list data (Price / Volume(1h) columns) per latest scrape.
# Note: some Price values are plain decimals (e.g. 313.27) which Excel would
# otherwise auto-convert to numbers; a leading apostrophe forces them to stay
# as text, matching the original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.541.30'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '2.492.30'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('E4').Value = '  -0.50%  '
$ws.Range('D5').Value = '''313.27'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').Value = '''94.10'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').Value = '  -0.93%  '
$ws.Range('D8').Value = '''0.998'
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('D9').Value = '''0.498'
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('D10').Value = '''32.84'
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E12').Value = '  +1.68%  '
$ws.Range('D13').Value = '2.876.61'
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').Value = '''6.86'
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').Value = '''15.50'
$ws.Range('E15').Value = '  +7.98%  '
$ws.Range('D16').Value = '2.469.02'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('E17').Value = '  -3.15%  '
$ws.Range('D18').Value = '41.613.53'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').Value = '''6.31'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').Value = '''70.75'
$ws.Range('E21').Value = '  +3.39%  '
$ws.Range('E22').Value = '  -2.38%  '
$ws.Range('D23').Value = '''236.01'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -1.16%  '
$ws.Range('D27').Value = '''24.56'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('E28').Value = '  +1.26%  '
$ws.Range('D29').Value = '''9.64'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').Value = '''36.32'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('D31').Value = '''154.10'
$ws.Range('E31').Value = '  +1.00%  '
$ws.Range('E32').Value = '  -2.83%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '''2.57'
$ws.Range('E33').Value = '  -2.34%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = '''18.16'
$ws.Range('E34').Value = '  +6.44%  '
$ws.Range('D35').Value = '''0.0758'
$ws.Range('E35').Value = '  +1.10%  '
$ws.Range('D36').Value = '''2.49'
$ws.Range('E36').Value = '  -2.05%  '
$ws.Range('E37').Value = '  -0.94%  '
$ws.Range('D38').Value = '''1.83'
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('E40').Value = '  -1.90%  '
$ws.Range('D41').Value = '''4.09'
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.44%  '
$ws.Range('D43').Value = '''19.61'
$ws.Range('E43').Value = '  -8.08%  '
$ws.Range('D44').Value = '1.950.13'
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('D45').Value = '''0.0284'
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('E46').Value = '  -2.12%  '
$ws.Range('D47').Value = '''8.81'
$ws.Range('E47').Value = '  +0.63%  '
$ws.Range('D48').Value = '2.729.93'
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('D49').Value = '''96.25'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('D51').Value = '''67.24'
$ws.Range('E51').Value = '  -2.10%  '
